$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the TPM-derived NATMI ligand-receptor metrics (rows 2-10) with
# recomputed values from the new TPM script run.
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1349983333333333
$ws.Range("H2").Value = 0.404995
$ws.Range("I2").Value = 0.06188478316908706
$ws.Range("J2").Value = 0.06188478316908706
$ws.Range("M2").Value = 1.599392
$ws.Range("N2").Value = 4.798176
$ws.Range("O2").Value = 0.03952976301548796
$ws.Range("P2").Value = 0.03952976301548796
$ws.Range("Q2").Value = 0.2159152543466666
$ws.Range("R2").Value = 1.94323728912
$ws.Range("S2").Value = 0.002446290812938869
$ws.Range("T2").Value = 0.002446290812938869
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1349983333333333
$ws.Range("H3").Value = 0.404995
$ws.Range("I3").Value = 0.06188478316908706
$ws.Range("J3").Value = 0.06188478316908706
$ws.Range("O3").Value = 0.4638329693976876
$ws.Range("P3").Value = 0.4638329693976876
$ws.Range("Q3").Value = 2.533498961848889
$ws.Range("R3").Value = 22.80149065664
$ws.Range("S3").Value = 0.02870420273784969
$ws.Range("T3").Value = 0.02870420273784969
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1349983333333333
$ws.Range("H4").Value = 0.404995
$ws.Range("I4").Value = 0.06188478316908706
$ws.Range("J4").Value = 0.06188478316908706
$ws.Range("M4").Value = 20.09416733333333
$ws.Range("N4").Value = 60.28250199999999
$ws.Range("O4").Value = 0.4966372675868244
$ws.Range("P4").Value = 0.4966372675868245
$ws.Range("Q4").Value = 2.712679099721111
$ws.Range("R4").Value = 24.41411189749
$ws.Range("S4").Value = 0.0307342896182985
$ws.Range("T4").Value = 0.0307342896182985
$ws.Range("I5").Value = 0.4284959871424753
$ws.Range("J5").Value = 0.4284959871424753
$ws.Range("M5").Value = 1.599392
$ws.Range("N5").Value = 4.798176
$ws.Range("O5").Value = 0.03952976301548796
$ws.Range("P5").Value = 0.03952976301548796
$ws.Range("Q5").Value = 1.495017277472
$ws.Range("R5").Value = 13.455155497248
$ws.Range("S5").Value = 0.01693834482482963
$ws.Range("T5").Value = 0.01693834482482963
$ws.Range("I6").Value = 0.4284959871424753
$ws.Range("J6").Value = 0.4284959871424753
$ws.Range("O6").Value = 0.4638329693976876
$ws.Range("P6").Value = 0.4638329693976876
$ws.Range("S6").Value = 0.1987505660912877
$ws.Range("T6").Value = 0.1987505660912877
$ws.Range("I7").Value = 0.4284959871424753
$ws.Range("J7").Value = 0.4284959871424753
$ws.Range("M7").Value = 20.09416733333333
$ws.Range("N7").Value = 60.28250199999999
$ws.Range("O7").Value = 0.4966372675868244
$ws.Range("P7").Value = 0.4966372675868245
$ws.Range("Q7").Value = 18.78284206732733
$ws.Range("R7").Value = 169.045578605946
$ws.Range("S7").Value = 0.212807076226358
$ws.Range("T7").Value = 0.212807076226358
$ws.Range("G8").Value = 1.111707
$ws.Range("H8").Value = 3.335121
$ws.Range("I8").Value = 0.5096192296884376
$ws.Range("J8").Value = 0.5096192296884376
$ws.Range("M8").Value = 1.599392
$ws.Range("N8").Value = 4.798176
$ws.Range("O8").Value = 0.03952976301548796
$ws.Range("P8").Value = 0.03952976301548796
$ws.Range("Q8").Value = 1.778055282144
$ws.Range("R8").Value = 16.002497539296
$ws.Range("S8").Value = 0.02014512737771947
$ws.Range("T8").Value = 0.02014512737771947
$ws.Range("G9").Value = 1.111707
$ws.Range("H9").Value = 3.335121
$ws.Range("I9").Value = 0.5096192296884376
$ws.Range("J9").Value = 0.5096192296884376
$ws.Range("O9").Value = 0.4638329693976876
$ws.Range("P9").Value = 0.4638329693976876
$ws.Range("Q9").Value = 20.863283722368
$ws.Range("R9").Value = 187.769553501312
$ws.Range("S9").Value = 0.2363782005685502
$ws.Range("T9").Value = 0.2363782005685502
$ws.Range("G10").Value = 1.111707
$ws.Range("H10").Value = 3.335121
$ws.Range("I10").Value = 0.5096192296884376
$ws.Range("J10").Value = 0.5096192296884376
$ws.Range("M10").Value = 20.09416733333333
$ws.Range("N10").Value = 60.28250199999999
$ws.Range("O10").Value = 0.4966372675868244
$ws.Range("P10").Value = 0.4966372675868245
$ws.Range("Q10").Value = 22.338826483638
$ws.Range("R10").Value = 201.049438352742
$ws.Range("S10").Value = 0.253095901742168
$ws.Range("T10").Value = 0.253095901742168
